$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Fill the previously-empty row 17 (between the existing row 16 and the
# trailing "Hoelang?"/footer rows 18-19) with a new log entry, copying the
# formatting/style banding from an existing "highlighted" row (row 13) so it
# matches the rest of the log (fill + border + wrap-text alignment).
$ws.Range("B13:D13").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B17").Value = "Udemy tutorials bekeken"
$ws.Range("C17").Value = "11/27/2021"
$ws.Range("D17").Value = "15 minuten"

$ws.Range("D23").Select()

$wb.Save()
